$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Break the merges whose regions need to move (sum row + footer row),
#    so the upcoming bulk value/format writes aren't blocked by a non-anchor
#    merged cell.
$ws.Range("K51:N51").UnMerge()
$ws.Range("A52:E52").UnMerge()
$ws.Range("F52:G52").UnMerge()
$ws.Range("I52:N52").UnMerge()

# 2) Shift formatting for the trailing special rows down by one row:
#    old footer (52) -> 53, old sum row (51) -> 52, and give the
#    soon-to-be-new data row (51) the same look as the other item rows.
$ws.Range("A52:N52").Copy()
$ws.Range("A53:N53").PasteSpecial(-4122)

$ws.Range("A51:N51").Copy()
$ws.Range("A52:N52").PasteSpecial(-4122)

$ws.Range("A50:N50").Copy()
$ws.Range("A51:N51").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3) Shift all the cell values for rows 34-52 down to rows 35-53 in one go.
$vals = $ws.Range("A34:N52").Value2
$ws.Range("A35:N53").Value2 = $vals

# 4) Row heights: new item row, shifted sum row, shifted footer row.
$ws.Rows(51).RowHeight = 24.75
$ws.Rows(52).RowHeight = 26.25
$ws.Rows(53).RowHeight = 16.5

# 5) Fill in the brand-new item row (34). A34 and N34 already carry the
#    right values from before the shift (31 and "1:0"), so only the name,
#    transaction string and quantity need to be set.
$ws.Range("B34").Value2 = "VOLTAREN 1% EMULGEL 25 GM"
$ws.Range("H34").Value2 = "1:0"
$ws.Range("L34").Value2 = 39

# 6) The total in K52 (old K51) needs to include the new row's quantity.
$ws.Range("K52").Value2 = 2789.8299999999999 + 39

# 7) Re-create merges: the new data row, and the relocated summary/footer rows.
$ws.Range("B51:G51").Merge()
$ws.Range("H51:K51").Merge()
$ws.Range("L51:M51").Merge()
$ws.Range("K52:N52").Merge()
$ws.Range("A53:E53").Merge()
$ws.Range("F53:G53").Merge()
$ws.Range("I53:N53").Merge()
